$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("8s window")

# Row-by-row data updates: B/C/D/E score columns, F (Depth) and G (Estimators)
$rows = @{
  5  = @{ B = 0.82383419689119097; C = 0.82383419689119097; D = 0.81102705277875597; E = 0.82562283797952796; F = 15;     G = 275 }
  6  = @{ B = 0.95535714285714202; C = 0.95535714285714202; D = 0.94081310960702003; E = 0.95542722583538897; F = 15;     G = 250 }
  7  = @{ B = 0.87309644670050701; C = 0.87309644670050701; D = 0.85661779081133904; E = 0.87367043741695705; F = 20;     G = 200 }
  8  = @{ B = 0.8125;              C = 0.8125;              D = 0.79407366209368402; E = 0.81528400576257498; F = 15;     G = 225 }
  9  = @{ B = 0.88695652173912998; C = 0.88695652173912998; D = 0.80882722857529799; E = 0.88328752556443602; F = 20;     G = 200 }
  10 = @{ B = 0.9;                 C = 0.9;                 D = 0.87087351556189496; E = 0.89705814960142005; F = "None"; G = 275 }
  11 = @{ B = 0.90640394088669896; C = 0.90640394088669896; D = 0.85485214950697197; E = 0.90302617691227505; F = 20;     G = 250 }
  12 = @{ B = 0.92307692307692302; C = 0.92307692307692302; D = 0.90966315148035504; E = 0.92284961247957398; F = "None"; G = 300 }
  13 = @{ B = 0.89839572192513295; C = 0.89839572192513295; D = 0.89719412137376497; E = 0.89961161250624599; F = 20;     G = 200 }
  14 = @{ B = 0.865979381443299;   C = 0.865979381443299;   D = 0.86880711880711803; E = 0.862081541978449;   F = 20;     G = 300 }
  15 = @{ B = 0.88888888888888795; C = 0.88888888888888795; D = 0.87948075893653699; E = 0.88848062162362595; F = "None"; G = 250 }
  16 = @{ B = 0.86528497409326399; C = 0.86528497409326399; D = 0.87236581048762496; E = 0.86932347679020905; F = 15;     G = 350 }
  17 = @{ B = 0.87562189054726303; C = 0.87562189054726303; D = 0.888942569868442;   E = 0.87479546587087298; F = "None"; G = 300 }
  18 = @{ B = 0.859375;            C = 0.859375;            D = 0.818063493456062;   E = 0.85897172992368398; F = "None"; G = 300 }
  19 = @{ B = 0.87755102040816302; C = 0.87755102040816302; D = 0.88119665527245095; E = 0.878369238553381;   F = 20;     G = 300 }
  20 = @{ B = 0.86274509803921495; C = 0.86274509803921495; D = 0.79089720218752402; E = 0.86298497695556498; F = 20;     G = 275 }
}

foreach ($r in $rows.Keys) {
  $vals = $rows[$r]
  $ws.Cells.Item($r, 2).Value = $vals.B
  $ws.Cells.Item($r, 3).Value = $vals.C
  $ws.Cells.Item($r, 4).Value = $vals.D
  $ws.Cells.Item($r, 5).Value = $vals.E
  if ($vals.F -eq "None") {
    $ws.Cells.Item($r, 6).Value = "None"
  } else {
    $ws.Cells.Item($r, 6).Value = $vals.F
  }
  $ws.Cells.Item($r, 7).Value = $vals.G
}

# Row 22 averages: now each of B..E has its own formula, F/G keep shared-formula look
$ws.Range("B22").Formula = "=SUM(B2:B21)/COUNT(B2:B21)"
$ws.Range("C22").Formula = "=SUM(C2:C21)/COUNT(C2:C21)"
$ws.Range("D22").Formula = "=SUM(D2:D21)/COUNT(D2:D21)"
$ws.Range("E22").Formula = "=SUM(E2:E21)/COUNT(E2:E21)"
$ws.Range("F22").Formula = "=SUM(F2:F21)/COUNT(F2:F21)"
$ws.Range("G22").Formula = "=SUM(G2:G21)/COUNT(G2:G21)"

# Selection state
$ws.Range("F26").Select()
